$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Author" -or $styleName -eq "Date") {
        $p.Alignment = 0
    }
}
